$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.547.43'
$ws.Range("E2").Value = '  +0.54%  '

$ws.Range("D3").Value = '3.475.52'
$ws.Range("E3").Value = '  +2.07%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '673.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.51'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.48%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.434'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("E9").Value = '  +1.64%  '

$ws.Range("E10").Value = '  +0.04%  '

$ws.Range("D11").Value = '3.473.58'
$ws.Range("E11").Value = '  +2.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +12.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.211'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.37%  '

$ws.Range("D14").Value = '98.568.23'
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.21'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000262'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("D17").Value = '4.130.46'
$ws.Range("E17").Value = '  +2.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.39%  '

$ws.Range("D19").Value = '3.471.51'
$ws.Range("E19").Value = '  +2.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.541'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '520.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.61%  '

$ws.Range("E24").Value = '  +0.94%  '

$ws.Range("E25").Value = '  +0.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '98.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.99%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.96%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +12.38%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.147'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.73%  '

$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.65%  '

$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.193'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.98%  '

$ws.Range("B34").Value = 'PolygonEcosystemToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.585'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.52%  '

$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.40%  '

$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '30.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.50%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.69%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.54'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.76%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '536.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.35%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.156'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.59%  '

$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.890'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.09%  '

$ws.Range("B43").Value = 'ImmutableX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.80%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0441'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.60%  '

$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '

$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.40%  '

$ws.Range("B48").Value = 'MantraDAO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.53%  '

$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.58%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '56.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.12%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.57%  '
